# Auto-generated edit script for LOM3256.docx restructuring
# Moves section headings' content blocks to the NEXT section (cyclic),
# and rotates the Avaliacao answers while pushing the Bibliografia text
# into the 'Norma de recuperacao' slot.
$d = $word.ActiveDocument

function Replace-InParagraph {
    param([int]$ParaIndex, [string]$OldText, [string]$NewText)
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $p.Range
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        Write-Output ('NOT FOUND in paragraph ' + $ParaIndex + ': ' + $OldText.Substring(0, [Math]::Min(60, $OldText.Length)))
    }
}

# Objetivos(PT) slot <- resumido(PT) text (paragraph 6)
Replace-InParagraph 6 'Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.^lO principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade^l(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estará apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superfícies de Fermi e constantes elásticas, usando um ou mais dos métodos e códigos computacionais apresentados em aula.' 'Revisão de mecânica quântica; Revisão de física do estado sólido; Método de Hartree-Fock; Teoria do funcional da densidade; Métodos de ondas planas e pseudo-potenciais; Códigos computacionais'

# Objetivos(EN) slot <- resumido(EN) text (paragraph 7)
Replace-InParagraph 7 'Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class.' 'Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes'

# Docente slot <- Objetivos(PT) text (paragraph 9)
Replace-InParagraph 9 '1176388 - Luiz Tadeu Fernandes Eleno' 'Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.^lO principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade^l(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estará apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superfícies de Fermi e constantes elásticas, usando um ou mais dos métodos e códigos computacionais apresentados em aula.'

# resumido(PT) slot <- Programa(PT) big text (paragraph 11)
Replace-InParagraph 11 'Revisão de mecânica quântica; Revisão de física do estado sólido; Método de Hartree-Fock; Teoria do funcional da densidade; Métodos de ondas planas e pseudo-potenciais; Códigos computacionais' 'Revisão de mecânica quântica^lo Equação de Schrödinger^lo Átomo do hidrogênio e orbitais atômicos^lo Notação de Dirac^lo Princípio variacional^lo Combinação linear de orbitais atômicos^lRevisão de física do estado sólido^lo Espaço direto e recíproco^lo Teorema de Bloch^lo Zona de Brillouin^lo Bandas de energia e densidade de estados^lo Energia de Fermi e superfície de Fermi^lo Aproximação de elétrons livres^lMétodo de Hartree-Fock^lo Determinantes de Slater^lo Equação de Hartree-Fock^lo Potencial de troca e correlação^lo Algoritmo autoconsistente^lTeoria do funcional da densidade^lo Teoremas de Hohenberg-Kohn^lo Equações de Kohn-Sham^lo Funcionais de troca e correlação: LDA, GGA, etc.^lMétodos de ondas planas e pseudo-potenciais^lo Bases de ondas planas^lo Pseudo-potenciais^lo Bases de ondas planas aumentadas e linearizadas^lo Método FP-LAPW^lCódigos computacionais^lo Quantum Espresso^lo Elk^lo Wien2k^lo VASP'

# resumido(EN) slot <- Objetivos(EN) text (paragraph 12)
Replace-InParagraph 12 'Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes' 'Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class.'

# Programa(PT) slot <- old Metodo answer text (paragraph 14)
Replace-InParagraph 14 'Revisão de mecânica quântica^lo Equação de Schrödinger^lo Átomo do hidrogênio e orbitais atômicos^lo Notação de Dirac^lo Princípio variacional^lo Combinação linear de orbitais atômicos^lRevisão de física do estado sólido^lo Espaço direto e recíproco^lo Teorema de Bloch^lo Zona de Brillouin^lo Bandas de energia e densidade de estados^lo Energia de Fermi e superfície de Fermi^lo Aproximação de elétrons livres^lMétodo de Hartree-Fock^lo Determinantes de Slater^lo Equação de Hartree-Fock^lo Potencial de troca e correlação^lo Algoritmo autoconsistente^lTeoria do funcional da densidade^lo Teoremas de Hohenberg-Kohn^lo Equações de Kohn-Sham^lo Funcionais de troca e correlação: LDA, GGA, etc.^lMétodos de ondas planas e pseudo-potenciais^lo Bases de ondas planas^lo Pseudo-potenciais^lo Bases de ondas planas aumentadas e linearizadas^lo Método FP-LAPW^lCódigos computacionais^lo Quantum Espresso^lo Elk^lo Wien2k^lo VASP' 'Aulas expositivas, trabalhos e exercícios comentados.'

# Bibliografia slot <- Docente text (paragraph 19)
Replace-InParagraph 19 'GRIFFITHS, D. J., Mecânica Quântica, Pearson.^lASHCROFT, N. W. Solid State Physics, Saunders College.^lKITTEL, C. Introduction to Solid State Physics. John Wiley & Sons.^lSUTTON, A. P. Electronic Structure of Materials, Oxford.^lMORGON, N. H. e COUTINHO, K. (eds), Métodos de Química teórica e modelagem molecular, Livraria da Física^lEditora.^lVIANNA, J. D. M., FAZZIO, A., CANUTO, S., Teoria Quântica de moléculas e sólidos, Livraria da Física Editora.^lCOTTENIER, S. Density Functional Theory and the Family of (L)APW-methods: a step-by-step introduction^l(apostila, disponível online)^lTHIJSSEN, J. M. Computational Physics, Cambridge.^lTADMOR, E. B., MILLER, R. E. Modeling Materials  Continuum, atomistic and multiscale techniques,^lCambridge.' '1176388 - Luiz Tadeu Fernandes Eleno'

# Avaliacao paragraph (17): Norma gets bibliography text,
# Criterio gets old Norma answer, Metodo gets old Criterio answer.
# (Order matters: process Norma, then Criterio, then Metodo so a
# freshly-written value is never re-matched by a later step.)
Replace-InParagraph 17 'Não haverá exame de recuperação' 'GRIFFITHS, D. J., Mecânica Quântica, Pearson.^lASHCROFT, N. W. Solid State Physics, Saunders College.^lKITTEL, C. Introduction to Solid State Physics. John Wiley & Sons.^lSUTTON, A. P. Electronic Structure of Materials, Oxford.^lMORGON, N. H. e COUTINHO, K. (eds), Métodos de Química teórica e modelagem molecular, Livraria da Física^lEditora.^lVIANNA, J. D. M., FAZZIO, A., CANUTO, S., Teoria Quântica de moléculas e sólidos, Livraria da Física Editora.^lCOTTENIER, S. Density Functional Theory and the Family of (L)APW-methods: a step-by-step introduction^l(apostila, disponível online)^lTHIJSSEN, J. M. Computational Physics, Cambridge.^lTADMOR, E. B., MILLER, R. E. Modeling Materials  Continuum, atomistic and multiscale techniques,^lCambridge.'
Replace-InParagraph 17 'Média aritmética de trabalhos propostos ao longo do curso.' 'Não haverá exame de recuperação'
Replace-InParagraph 17 'Aulas expositivas, trabalhos e exercícios comentados.' 'Média aritmética de trabalhos propostos ao longo do curso.'
